$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updates to nombre_aides (column C) and montant_total (column E) for rows
# affected by the 2022-05-09 data refresh. nombre_entreprises (column D) is
# unchanged.
$updates = @(
    @{ Row = 3;   C = 249328; E = 1036478313 },
    @{ Row = 53;  C = 141680; E = 590063206 },
    @{ Row = 57;  C = 3712;   E = 138408892 },
    @{ Row = 92;  C = 409074; E = 1594755556 },
    @{ Row = 93;  C = 209550; E = 1308407897 },
    @{ Row = 95;  C = 50755;  E = 931923076 },
    @{ Row = 96;  C = 17257;  E = 790651633 },
    @{ Row = 104; C = 135233; E = 272168047 },
    @{ Row = 110; C = 396;    E = 16649846 },
    @{ Row = 174; C = 226089; E = 900577214 },
    @{ Row = 175; C = 80780;  E = 486154029 }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 3).Value = $u.C
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}
